# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet with the latest scraped figures, matching the GitHub Actions
# scraper commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some D-column prices are plain decimal-looking strings (e.g. "12.80").
# Forcing the cell to Text format ("@") before assigning the value keeps
# Excel from reinterpreting them as numbers and silently dropping
# significant trailing zeros, matching the worksheets original
# inline-string (text) cell type.

$ws.Range('D2').Value = '35.325.61'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').Value = '1.912.48'
$ws.Range('E3').Value = '  +0.24%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.723'
$ws.Range('E5').Value = '  +8.99%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '254.77'
$ws.Range('E6').Value = '  +3.67%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '40.76'
$ws.Range('E8').Value = '  -1.46%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.370'
$ws.Range('E9').Value = '  +5.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.85'
$ws.Range('E10').Value = '  +0.09%  '
$ws.Range('E11').Value = '  +6.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0985'
$ws.Range('E12').Value = '  -0.81%  '
$ws.Range('D13').Value = '2.187.00'
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.80'
$ws.Range('E14').Value = '  +5.87%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.725'
$ws.Range('E15').Value = '  +3.66%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.97'
$ws.Range('E16').Value = '  +2.21%  '
$ws.Range('D17').Value = '1.911.88'
$ws.Range('E17').Value = '  -0.14%  '
$ws.Range('D18').Value = '35.330.34'
$ws.Range('E18').Value = '  -0.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '74.50'
$ws.Range('E19').Value = '  +2.66%  '
$ws.Range('D20').Value = '0.0₃0855'
$ws.Range('E20').Value = '  +3.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '244.07'
$ws.Range('E21').Value = '  +1.90%  '
$ws.Range('E22').Value = '  +4.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.12'
$ws.Range('E23').Value = '  +6.03%  '
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.47'
$ws.Range('E25').Value = '  +5.35%  '
$ws.Range('E26').Value = '  +4.04%  '
$ws.Range('E27').Value = '  -1.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.66'
$ws.Range('E29').Value = '  +1.89%  '
$ws.Range('E30').Value = '  +4.93%  '
$ws.Range('D31').Value = '4.129.63'
$ws.Range('E31').Value = '  +19.48%  '
$ws.Range('E32').Value = '  +5.28%  '
$ws.Range('E33').Value = '  +14.16%  '
$ws.Range('E34').Value = '  +23.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0588'
$ws.Range('E35').Value = '  +4.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.23'
$ws.Range('E36').Value = '  +2.95%  '
$ws.Range('E37').Value = '  -0.15%  '
$ws.Range('E38').Value = '  -2.71%  '
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0219'
$ws.Range('E40').Value = '  +5.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '17.13'
$ws.Range('E41').Value = '  +5.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '96.96'
$ws.Range('E42').Value = '  +7.86%  '
$ws.Range('E43').Value = '  +1.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0645'
$ws.Range('E44').Value = '  +1.23%  '
$ws.Range('D45').Value = '1.336.64'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.43'
$ws.Range('E46').Value = '  +1.76%  '
$ws.Range('E47').Value = '  +0.97%  '
$ws.Range('E48').Value = '  +3.00%  '
$ws.Range('E49').Value = '  -0.62%  '
$ws.Range('E50').Value = '  -5.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.87'
$ws.Range('E51').Value = '  +12.28%  '
